$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("C1").Value = "Resolución Primigenia"
$ws.Range("D1").Value = "Resolución Hija"

# Clear the "Resolución Hija" column values for the sample rows
$ws.Range("D2").Value = ""
$ws.Range("D3").Value = ""
